$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Simple text edits (Find/Replace keeps the existing run/paragraph
#    formatting intact, it only swaps the literal text of the run).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Yhteys Puttylla", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Yhteys työasemasta Puttylla", 2) | Out-Null

$d.Content.Find.Execute("Staattinen IP (ongelmia?)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Staattinen IP", 2) | Out-Null

$d.Content.Find.Execute("Linux Xubuntu VM (Database)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Linux Xubuntu VM (Tietokanta)", 2) | Out-Null

$d.Content.Find.Execute("Ei saaa yhteyttä mihinkään (ping)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Tietokanta luotu", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Structural insertions.
#
# The engine's Range.InsertXML grafts the *last* <w:p> of the inserted
# fragment onto the paragraph mark of whatever already follows the
# insertion point, so every fragment below ends with a "spacer" paragraph
# whose formatting mirrors that following paragraph - this keeps the
# graft invisible while still landing all the genuinely new paragraphs
# in between with exactly the formatting we want.
# ---------------------------------------------------------------------------

function New-OpenXmlPackage($bodyXml) {
    return @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
}

# --- Insertion point #1: right after "Staattinen IP" ----------------------
# New paragraphs needed, in order:
#   "SSH yhteys Tietokantaan"                  (not bold)
#   ""                                          (not bold)
#   "Python koodi"                              (bold)
#   "Hakee raspi dataa ja lisää sen kantaan"    (not bold)
# followed by the pre-existing empty bold paragraph (kept as-is via the
# spacer trick below).

$pStaattinen = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Staattinen IP`r") {
        $pStaattinen = $d.Paragraphs.Item($i)
        break
    }
}

$insertPoint1 = $d.Range($pStaattinen.Range.End, $pStaattinen.Range.End)
$body1 = @'
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>SSH yhteys Tietokantaan</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Python koodi</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>Hakee raspi dataa ja lisää sen kantaan</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr></w:r></w:p>
'@
$insertPoint1.InsertXML((New-OpenXmlPackage $body1))

# --- Insertion point #2: right after "Tietokanta luotu" --------------------
# New paragraph needed:
#   ""                                          (not bold)
# followed by the pre-existing empty bold paragraph (kept as-is via the
# spacer trick below).

$pTietokantaLuotu = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Tietokanta luotu`r") {
        $pTietokantaLuotu = $d.Paragraphs.Item($i)
        break
    }
}

$insertPoint2 = $d.Range($pTietokantaLuotu.Range.End, $pTietokantaLuotu.Range.End)
$body2 = @'
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr></w:r></w:p>
'@
$insertPoint2.InsertXML((New-OpenXmlPackage $body2))
